$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 42-44 in column A were plain RUT numbers (or an inconsistent shared
# string) without the dash separator used elsewhere in the sheet. Normalize
# them to the "NNNNNNNN-D" text format, which also means they become shared
# strings instead of numeric values.
$ws.Range("A42").Value = "16610707-5"
$ws.Range("A43").Value = "18462110-K"
$ws.Range("A44").Value = "18741199-8"

# Move the active selection down to the last data row (A44).
[void]$ws.Range("A44").Select()
